# Update "想去人数" (interest count, column F) values on the 展览, 演出 and
# 全部类型 sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1104
$ws1.Range("F12").Value = 563
$ws1.Range("F14").Value = 1783
$ws1.Range("F15").Value = 830
$ws1.Range("F17").Value = 1448
$ws1.Range("F21").Value = 388
$ws1.Range("F24").Value = 4657
$ws1.Range("F25").Value = 730
$ws1.Range("F27").Value = 1614
$ws1.Range("F29").Value = 86

# --- 演出 (Performances) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 49
$ws2.Range("F14").Value = 22

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 49
$ws4.Range("F15").Value = 1104
$ws4.Range("F23").Value = 563
$ws4.Range("F25").Value = 1783
$ws4.Range("F26").Value = 830
$ws4.Range("F28").Value = 1448
$ws4.Range("F34").Value = 388
$ws4.Range("F37").Value = 4657
$ws4.Range("F38").Value = 730
$ws4.Range("F40").Value = 1614
$ws4.Range("F42").Value = 22
$ws4.Range("F44").Value = 86
